$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.330.66"
$ws.Range("E2").Value = "  +0.41%  "

# Row 3
$ws.Range("D3").Value = "2.238.78"
$ws.Range("E3").Value = "  -0.83%  "

# Row 4
$ws.Range("E4").Value = "  +0.51%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.52"
$ws.Range("E5").Value = "  -0.94%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.33"
$ws.Range("E6").Value = "  -6.73%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("E7").Value = "  -0.75%  "

# Row 8
$ws.Range("E8").Value = "  +0.38%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("E9").Value = "  -3.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.36"
$ws.Range("E10").Value = "  -3.67%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("E11").Value = "  -1.93%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.11"
$ws.Range("E12").Value = "  -3.36%  "

# Row 13
$ws.Range("E13").Value = "  -0.03%  "

# Row 14
$ws.Range("D14").Value = "2.363.37"
$ws.Range("E14").Value = "  +4.46%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.831"
$ws.Range("E15").Value = "  -1.33%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.49"
$ws.Range("E16").Value = "  -3.04%  "

# Row 17
$ws.Range("D17").Value = "44.008.63"
$ws.Range("E17").Value = "  -0.10%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0958"
$ws.Range("E18").Value = "  -1.80%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.30"
$ws.Range("E19").Value = "  -4.62%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.32"
$ws.Range("E20").Value = "  -0.70%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "65.49"
$ws.Range("E21").Value = "  -0.04%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.12"
$ws.Range("E22").Value = "  +5.17%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.46"
$ws.Range("E23").Value = "  -1.86%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("E24").Value = "  -0.68%  "

# Row 25
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.57"
$ws.Range("E26").Value = "  +2.41%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.20"
$ws.Range("E27").Value = "  +5.73%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.74"
$ws.Range("E28").Value = "  -4.45%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.91"
$ws.Range("E29").Value = "  -4.55%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.94"
$ws.Range("E30").Value = "  -1.04%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.36"
$ws.Range("E31").Value = "  -2.56%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0793"
$ws.Range("E32").Value = "  -3.75%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.63"
$ws.Range("E33").Value = "  -1.22%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.09"
$ws.Range("E34").Value = "  -13.32%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.108"
$ws.Range("E35").Value = "  +0.56%  "

# Row 36
$ws.Range("E36").Value = "  -0.47%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.76"
$ws.Range("E37").Value = "  -5.71%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.44"
$ws.Range("E38").Value = "  +1.15%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.59"
$ws.Range("E39").Value = "  -7.80%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.79"
$ws.Range("E40").Value = "  -3.14%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0299"
$ws.Range("E41").Value = "  -2.42%  "

# Row 42
$ws.Range("E42").Value = "  +0.38%  "

# Row 43
$ws.Range("D43").Value = "1.732.75"
$ws.Range("E43").Value = "  -2.56%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "80.04"
$ws.Range("E44").Value = "  -9.04%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.190"
$ws.Range("E45").Value = "  -1.91%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "99.06"
$ws.Range("E46").Value = "  -2.75%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +3.98%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.91"
$ws.Range("E48").Value = "  -4.89%  "

# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.12"
$ws.Range("E49").Value = "  -2.15%  "

# Row 50
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.92"
$ws.Range("E50").Value = "  -1.25%  "

# Row 51
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.96"
$ws.Range("E51").Value = "  -3.46%  "
